$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.854.59'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '1.936.47'
$ws.Range('E3').Value = '  -1.03%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.09%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '243.54'
$c.ClearFormats()
$ws.Range('E6').Value = '  -0.10%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4902'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.29%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2957'
$c.ClearFormats()
$ws.Range('E8').Value = '  -0.54%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06899'
$c.ClearFormats()
$ws.Range('E9').Value = '  +0.65%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '19.31'
$c.ClearFormats()
$ws.Range('E10').Value = '  +0.40%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '104.90'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.84%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.07796'
$c.ClearFormats()
$ws.Range('E12').Value = '  +0.46%  '
$ws.Range('D13').Value = '1.934.64'
$ws.Range('E13').Value = '  -1.05%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '5.352'
$c.ClearFormats()
$ws.Range('E14').Value = '  -2.12%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.7025'
$c.ClearFormats()
$ws.Range('E15').Value = '  -0.78%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '273.41'
$c.ClearFormats()
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('D17').Value = '30.844.55'
$ws.Range('E17').Value = '  -0.79%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000007735'
$c.ClearFormats()
$ws.Range('E18').Value = '  -0.35%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '13.08'
$c.ClearFormats()
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('E20').Value = '  -0.07%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.579'
$c.ClearFormats()
$ws.Range('E21').Value = '  +1.26%  '
$ws.Range('D22').Value = '2.187.59'
$ws.Range('E22').Value = '  -0.39%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('E24').Value = '  +0.29%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.863'
$c.ClearFormats()
$ws.Range('E25').Value = '  +0.25%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '166.26'
$c.ClearFormats()
$ws.Range('E26').Value = '  -2.12%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '19.55'
$c.ClearFormats()
$ws.Range('E27').Value = '  -2.57%  '
$ws.Range('E28').Value = '  -2.71%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.1042'
$c.ClearFormats()
$ws.Range('E29').Value = '  -1.34%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '1.389'
$c.ClearFormats()
$ws.Range('E30').Value = '  -2.91%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '1.562'
$c.ClearFormats()
$ws.Range('E31').Value = '  -1.56%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '4.568'
$c.ClearFormats()
$ws.Range('E32').Value = '  -0.29%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.377'
$c.ClearFormats()
$ws.Range('E33').Value = '  -1.88%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.04887'
$c.ClearFormats()
$ws.Range('E34').Value = '  -1.98%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.7626'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.40%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.150'
$c.ClearFormats()
$ws.Range('E36').Value = '  -2.98%  '
$ws.Range('E37').Value = '  -0.09%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.733'
$c.ClearFormats()
$ws.Range('E38').Value = '  -0.07%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.02009'
$c.ClearFormats()
$ws.Range('E39').Value = '  -1.45%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '79.00'
$c.ClearFormats()
$ws.Range('E40').Value = '  +4.70%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.656'
$c.ClearFormats()
$ws.Range('E41').Value = '  -1.89%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '6.513'
$c.ClearFormats()
$ws.Range('E42').Value = '  +0.08%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '2.089'
$c.ClearFormats()
$ws.Range('E43').Value = '  -4.12%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.9066'
$c.ClearFormats()
$ws.Range('E44').Value = '  +2.30%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.4441'
$c.ClearFormats()
$ws.Range('E45').Value = '  -1.66%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '107.91'
$c.ClearFormats()
$ws.Range('E46').Value = '  -1.32%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '7.726'
$c.ClearFormats()
$ws.Range('E48').Value = '  -4.84%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '994.79'
$c.ClearFormats()
$ws.Range('E49').Value = '  +1.03%  '
$ws.Range('E50').Value = '  -1.35%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '36.17'
$c.ClearFormats()
$ws.Range('E51').Value = '  +0.90%  '
